$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted as row 87, pushing the
# previously existing rows 87-95 down to 88-96 (dimension grows to R96).
$ws.Rows.Item(87).Insert()

# Populate the newly inserted row 87 with the new record's data.
$ws.Cells.Item(87, 1).Value = 10
$ws.Cells.Item(87, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(87, 3).Value = "La Araucanía"
$ws.Cells.Item(87, 4).Value = 45154
$ws.Cells.Item(87, 5).Value = 9
$ws.Cells.Item(87, 6).Value = 100112042
$ws.Cells.Item(87, 7).Value = "Locoto"
$ws.Cells.Item(87, 8).Value = "Sin especificar"
$ws.Cells.Item(87, 9).Value = "Primera"
$ws.Cells.Item(87, 10).Value = 80
$ws.Cells.Item(87, 11).Value = 2700
$ws.Cells.Item(87, 12).Value = 2700
$ws.Cells.Item(87, 13).Value = 2700
$ws.Cells.Item(87, 14).Value = "`$/kilo"
$ws.Cells.Item(87, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(87, 16).Value = 2700
$ws.Cells.Item(87, 17).Value = 1
$ws.Cells.Item(87, 18).Value = "Hortaliza"
